# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Update COVID figures for Estados Unidos (USA), and for Canada, Uzbekistan,
#   Ruanda and Siria, whose rows moved up (new/updated figures) while the
#   countries they displaced keep their previous figures in the following row(s).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Timestamp update (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 19:52"

# Estados Unidos (row 4) - figures refreshed in place
Set-CountryRow 4 "Estados Unidos" 622412 8526 47707 547156 13477 1502 27549

# Canada moves above Paises Bajos with refreshed figures (rows 14-15)
Set-CountryRow 14 "Canada" 28205 1142 8937 18262 557 103 1006
Set-CountryRow 15 "Paises Bajos" 28153 734 250 24769 1279 189 3134

# Uzbekistan moves above Kazajistan with refreshed figures (rows 69-70)
Set-CountryRow 69 "Uzbekistan" 1302 137 107 1191 8 0 4
Set-CountryRow 70 "Kazajistan" 1295 63 240 1039 20 2 16

# Ruanda moves above Brunei with refreshed figures (rows 128-129)
Set-CountryRow 128 "Ruanda" 136 2 54 82 0 0 0
Set-CountryRow 129 "Brunei" 136 0 108 27 2 0 1

# Siria moves above Guam/Sudan/Mongolia/Mozambique with refreshed figures (rows 166-170)
Set-CountryRow 166 "Siria" 33 4 5 26 0 0 2
Set-CountryRow 167 "Guam" 32 0 0 31 0 0 1
Set-CountryRow 168 "Sudan" 32 0 4 23 0 0 5
Set-CountryRow 169 "Mongolia" 30 0 5 25 0 0 0
Set-CountryRow 170 "Mozambique" 29 1 2 27 0 0 0
